$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 7117.5884
$ws.Range("I2").Value = 1339.0714
$ws.Range("K2").Value = 1339.0714
$ws.Range("M2").Value = -1226.0714
$ws.Range("H40").Value = 2332.923
$ws.Range("I40").Value = 2498.625
$ws.Range("J40").Value = 2067.8
$ws.Range("K40").Value = 2498.625
$ws.Range("L40").Value = 2067.8
$ws.Range("M40").Value = -2323.625
$ws.Range("N40").Value = -2417.8
$ws.Range("H53").Value = 909.8461
$ws.Range("I53").Value = 119.833336
$ws.Range("K53").Value = 119.833336
$ws.Range("M53").Value = 517.166664
$ws.Range("H62").Value = 38097740
$ws.Range("I62").Value = 44447028
$ws.Range("J62").Value = 2000
$ws.Range("K62").Value = 44447028
$ws.Range("L62").Value = 2000
$ws.Range("M62").Value = -44446404
$ws.Range("N62").Value = -3248
$ws.Range("H65").Value = 38097740
$ws.Range("I65").Value = 44447028
$ws.Range("J65").Value = 2000
$ws.Range("K65").Value = 222235140
$ws.Range("L65").Value = 10000
$ws.Range("M65").Value = -222232020
$ws.Range("N65").Value = -16240
$ws.Range("H88").Value = 8565
$ws.Range("I88").Value = 6173.75
$ws.Range("K88").Value = 6173.75
$ws.Range("M88").Value = -5767.75
$ws.Range("H91").Value = 8565
$ws.Range("I91").Value = 6173.75
$ws.Range("K91").Value = 6173.75
$ws.Range("M91").Value = -4769.75
$ws.Range("H100").Value = 8425.114
$ws.Range("I100").Value = 2175.8
$ws.Range("K100").Value = 2175.8
$ws.Range("M100").Value = -1634.8
$ws.Range("H106").Value = 2833
$ws.Range("I106").Value = 2833
$ws.Range("K106").Value = 2833
$ws.Range("M106").Value = -2202
$ws.Range("H137").Value = 20842308
$ws.Range("I137").Value = 50001096
$ws.Range("J137").Value = 14602.571
$ws.Range("K137").Value = 150003288
$ws.Range("L137").Value = 43807.713
$ws.Range("M137").Value = -150000738
$ws.Range("N137").Value = -48907.713

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21671.773
$ws.Range("I32").Value = 23218.348
$ws.Range("K32").Value = 23218.348
$ws.Range("M32").Value = -22931.348
$ws.Range("H41").Value = 3762.111
$ws.Range("I41").Value = 3224.625
$ws.Range("K41").Value = 3224.625
$ws.Range("M41").Value = -2810.625
$ws.Range("H50").Value = 414.125
$ws.Range("I50").Value = 477.66666
$ws.Range("J50").Value = 376
$ws.Range("K50").Value = 477.66666
$ws.Range("L50").Value = 376
$ws.Range("M50").Value = 236.33334
$ws.Range("N50").Value = -1804
$ws.Range("H112").Value = 18720
$ws.Range("J112").Value = 18720
$ws.Range("L112").Value = 18720
$ws.Range("N112").Value = -21674
$ws.Range("H122").Value = 2562.8333
$ws.Range("I122").Value = 1722
$ws.Range("K122").Value = 5166
$ws.Range("M122").Value = -2716

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 723.1111
$ws.Range("I12").Value = 118
$ws.Range("J12").Value = 1933.3334
$ws.Range("K12").Value = 118
$ws.Range("L12").Value = 1933.3334
$ws.Range("M12").Value = 50
$ws.Range("N12").Value = -2269.3334
$ws.Range("H20").Value = 2559.0688
$ws.Range("I20").Value = 2675
$ws.Range("J20").Value = 2434.8572
$ws.Range("K20").Value = 2675
$ws.Range("L20").Value = 2434.8572
$ws.Range("M20").Value = -2428
$ws.Range("N20").Value = -2928.8572
$ws.Range("H46").Value = 20890
$ws.Range("J46").Value = 20890
$ws.Range("L46").Value = 20890
$ws.Range("N46").Value = -21486
$ws.Range("H75").Value = 6552.6665
$ws.Range("I75").Value = 3646.75
$ws.Range("J75").Value = 29800
$ws.Range("K75").Value = 3646.75
$ws.Range("L75").Value = 29800
$ws.Range("M75").Value = -2710.75
$ws.Range("N75").Value = -31672
$ws.Range("H78").Value = 6552.6665
$ws.Range("I78").Value = 3646.75
$ws.Range("J78").Value = 29800
$ws.Range("K78").Value = 10940.25
$ws.Range("L78").Value = 89400
$ws.Range("M78").Value = -6260.25
$ws.Range("N78").Value = -98760

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 1880.1111
$ws.Range("I3").Value = 240.125
$ws.Range("K3").Value = 240.125
$ws.Range("M3").Value = -127.125
$ws.Range("H22").Value = 1351.381
$ws.Range("I22").Value = 418.92856
$ws.Range("K22").Value = 418.92856
$ws.Range("M22").Value = -68.92856
$ws.Range("H31").Value = 62505520
$ws.Range("I31").Value = 333334340
$ws.Range("J31").Value = 6566.154
$ws.Range("K31").Value = 333334340
$ws.Range("L31").Value = 6566.154
$ws.Range("M31").Value = -333334045
$ws.Range("N31").Value = -7156.154
$ws.Range("H34").Value = 62505520
$ws.Range("I34").Value = 333334340
$ws.Range("J34").Value = 6566.154
$ws.Range("K34").Value = 333334340
$ws.Range("L34").Value = 6566.154
$ws.Range("M34").Value = -333334138
$ws.Range("N34").Value = -6970.154
$ws.Range("H62").Value = 8608.777
$ws.Range("J62").Value = 9855
$ws.Range("L62").Value = 9855
$ws.Range("N62").Value = -11103
$ws.Range("H65").Value = 8608.777
$ws.Range("J65").Value = 9855
$ws.Range("L65").Value = 49275
$ws.Range("N65").Value = -55515
$ws.Range("H103").Value = 15502.75
$ws.Range("I103").Value = 15502.75
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 15502.75
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -14330.75
$ws.Range("N103").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 5527.4614
$ws.Range("J2").Value = 8965.125
$ws.Range("L2").Value = 53790.75
$ws.Range("N2").Value = -54016.75
$ws.Range("H4").Value = 30470814
$ws.Range("I4").Value = 42268110
$ws.Range("K4").Value = 126804330
$ws.Range("M4").Value = -126804218
$ws.Range("H15").Value = 3721.879
$ws.Range("J15").Value = 5125.8696
$ws.Range("L15").Value = 15377.6088
$ws.Range("N15").Value = -15657.6088
$ws.Range("H35").Value = 3843.1667
$ws.Range("I35").Value = 574.5
$ws.Range("K35").Value = 1723.5
$ws.Range("M35").Value = -1435.5
$ws.Range("H38").Value = 61.75
$ws.Range("I38").Value = 21.90909
$ws.Range("K38").Value = 65.72727
$ws.Range("M38").Value = 281.27273
$ws.Range("H120").Value = 10000
$ws.Range("I120").Value = 10000
$ws.Range("K120").Value = 30000
$ws.Range("M120").Value = -25162
$ws.Range("H132").Value = 1915
$ws.Range("J132").Value = 4490
$ws.Range("L132").Value = 40410
$ws.Range("N132").Value = -45470

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 45.666668
$ws.Range("I13").Value = 45.666668
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 45.666668
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 93.333332
$ws.Range("N13").ClearContents()
$ws.Range("H41").Value = 2324.5
$ws.Range("I41").Value = 971.7143
$ws.Range("J41").Value = 4218.4
$ws.Range("K41").Value = 971.7143
$ws.Range("L41").Value = 4218.4
$ws.Range("M41").Value = -616.7143
$ws.Range("N41").Value = -4928.4
$ws.Range("H122").Value = 6569.3
$ws.Range("I122").Value = 6937.952
$ws.Range("K122").Value = 20813.856
$ws.Range("M122").Value = -18363.856
$ws.Range("H132").Value = 12009.4
$ws.Range("J132").Value = 13849.5
$ws.Range("L132").Value = 41548.5
$ws.Range("N132").Value = -46608.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 9306
$ws.Range("I46").Value = 3500
$ws.Range("J46").Value = 10080.134
$ws.Range("K46").Value = 3500
$ws.Range("L46").Value = 10080.134
$ws.Range("M46").Value = -3312
$ws.Range("N46").Value = -10456.134
$ws.Range("H104").Value = 20000
$ws.Range("J104").Value = 20000
$ws.Range("L104").Value = 20000
$ws.Range("N104").Value = -26988
$ws.Range("H110").Value = 88559.336
$ws.Range("J110").Value = 88559.336
$ws.Range("L110").Value = 88559.336
$ws.Range("N110").Value = -96739.336
$ws.Range("H122").Value = 3493.2307
$ws.Range("I122").Value = 3240.6
$ws.Range("J122").Value = 4335.3335
$ws.Range("K122").Value = 9721.799999999999
$ws.Range("L122").Value = 13006.0005
$ws.Range("M122").Value = -7271.799999999999
$ws.Range("N122").Value = -17906.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 48999
$ws.Range("J47").Value = 48999
$ws.Range("L47").Value = 48999
$ws.Range("N47").Value = -50143
$ws.Range("H52").Value = 14970.833
$ws.Range("I52").Value = 6697.4443
$ws.Range("K52").Value = 6697.4443
$ws.Range("M52").Value = -6471.4443
$ws.Range("H122").Value = 4296.5
$ws.Range("I122").Value = 4390.5557
$ws.Range("K122").Value = 13171.6671
$ws.Range("M122").Value = -10721.6671
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H124").Value = 108474
$ws.Range("J124").Value = 108474
$ws.Range("L124").Value = 108474
$ws.Range("N124").Value = -118294
$ws.Range("H136").Value = 3294.923
$ws.Range("I136").Value = 2200.125
$ws.Range("K136").Value = 6600.375
$ws.Range("M136").Value = -4050.375
